$d = $word.ActiveDocument
$sep = "--------------------------------------------------------------------------------"

# --- Change 1: insert the new "Date 1 day" / "Date 2 day" blocks right before "Affidavit:" ---
$target = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    if ($d.Paragraphs.Item($i).Range.Text -eq "Affidavit:$([char]13)") {
        $target = $d.Paragraphs.Item($i)
        break
    }
}
if ($target -eq $null) {
    throw "Could not locate 'Affidavit:' paragraph"
}
$insertText = "Date 1 day:`r{{ DATE1DAY_NUMBER }}`r$sep`rDate 2 day:`r{{ DATE2DAY_NUMBER }}`r$sep`r"
$target.Range.InsertBefore($insertText)

# --- Change 2: mark the paragraph containing {{ NIGHTJUSTIFY }} with a lastRenderedPageBreak ---
$njPara = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    if ($d.Paragraphs.Item($i).Range.Text -eq "{{ NIGHTJUSTIFY }}$([char]13)") {
        $njPara = $d.Paragraphs.Item($i)
        break
    }
}
if ($njPara -eq $null) {
    throw "Could not locate '{{ NIGHTJUSTIFY }}' paragraph"
}
$njRange = $njPara.Range.Duplicate
$xml = '<?xml version="1.0" standalone="yes"?><?mso-application progid="Word.Document"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r><w:lastRenderedPageBreak/><w:t>{{ NIGHTJUSTIFY }}</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$njRange.InsertXML($xml) | Out-Null

Write-Output "edit complete"
